# 408: Update unit tests and excel files for outstanding reports
#
# Adds a new "WMT_Extract_SA" worksheet (WMT_Extract "Standalone Order"
# extract) immediately after the existing "T2A" sheet, and gives it the
# same header-row layout/style used by the other extract sheets plus three
# new trailing columns (Disposal_Type_Desc, Disposal_Type_Code,
# Standalone_Order).

$wb = $excel.ActiveWorkbook

# Insert the new sheet right after T2A (the last sheet) - it becomes the
# active / selected sheet, which is what makes it the new "tabSelected"
# sheet and bumps the workbook's ActiveTab.
$t2a = $wb.Worksheets.Item("T2A")
$ws = $wb.Worksheets.Add($null, $t2a)
$ws.Name = "WMT_Extract_SA"

# Borrow the header cell formatting (fill/border/font/number format) that
# the other extract sheets already use for their header row, so the shared
# style is reused rather than a new one being created.
$styleSource = $wb.Worksheets.Item("Court_Reports").Range("C1")
$styleSource.Copy()
$ws.Range("A1:I1").PasteSpecial(-4122)  # xlPasteFormats

# Header row values (shared strings are reused where they already exist
# elsewhere in the workbook, and appended when new).
$ws.Range("A1").Value = "Case_Ref_No"
$ws.Range("B1").Value = "Tier_Code"
$ws.Range("C1").Value = "Team_Code"
$ws.Range("D1").Value = "OM_Grade_Code"
$ws.Range("E1").Value = "OM_Key"
$ws.Range("F1").Value = "Location"
$ws.Range("G1").Value = "Disposal_Type_Desc"
$ws.Range("H1").Value = "Disposal_Type_Code"
$ws.Range("I1").Value = "Standalone_Order"

# Match the taller header row height used on the other extract sheets.
$ws.Rows.Item(1).RowHeight = 18

# Leave the same cell selected/active as in the authored workbook.
$ws.Range("L10").Select()
